$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.534.07'
$ws.Range("E2").Value = '  -2.72%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.860.46'
$ws.Range("E3").Value = '  -2.57%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '286.75'
$ws.Range("E5").Value = '  -6.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5214'
$ws.Range("E7").Value = '  -3.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3684'
$ws.Range("E8").Value = '  -3.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07079'
$ws.Range("E9").Value = '  -2.94%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.01'
$ws.Range("E10").Value = '  -4.79%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8740'
$ws.Range("E11").Value = '  -3.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08052'
$ws.Range("E12").Value = '  -1.79%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.854.24'
$ws.Range("E13").Value = '  +61.08%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.218'
$ws.Range("E14").Value = '  -2.40%  '

$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.49'
$ws.Range("E15").Value = '  -5.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  -0.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.52'
$ws.Range("E17").Value = '  -2.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008395'
$ws.Range("E18").Value = '  -3.00%  '

$ws.Range("E19").Value = '  +0.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.600.03'
$ws.Range("E20").Value = '  -2.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.912'
$ws.Range("E21").Value = '  -2.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.51'
$ws.Range("E22").Value = '  -2.85%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.301'
$ws.Range("E23").Value = '  -3.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.70'
$ws.Range("E24").Value = '  -3.53%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.221'
$ws.Range("E25").Value = '  -3.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.737'
$ws.Range("E26").Value = '  -0.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.76'
$ws.Range("E27").Value = '  -2.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '112.72'
$ws.Range("E28").Value = '  -3.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.623'
$ws.Range("E29").Value = '  -4.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.530'
$ws.Range("E30").Value = '  -5.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09001'
$ws.Range("E31").Value = '  -3.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7816'
$ws.Range("E32").Value = '  -6.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04916'
$ws.Range("E33").Value = '  -3.06%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.148'
$ws.Range("E34").Value = '  -6.41%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.910'
$ws.Range("E35").Value = '  -2.91%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.5811'
$ws.Range("E36").Value = '  +0.87%  '

$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.164'
$ws.Range("E37").Value = '  -5.74%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.568'
$ws.Range("E38").Value = '  -5.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01918'
$ws.Range("E39").Value = '  -4.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.047'
$ws.Range("E40").Value = '  -2.82%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.378'
$ws.Range("E41").Value = '  -2.80%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '114.08'
$ws.Range("E42").Value = '  -2.91%  '

$ws.Range("B43").Value = 'Decentraland'
$ws.Range("C43").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5085'
$ws.Range("E43").Value = '  +3.04%  '

$ws.Range("E44").Value = '  -8.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1465'
$ws.Range("E45").Value = '  -3.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.836'
$ws.Range("E47").Value = '  -2.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.593'
$ws.Range("E48").Value = '  -2.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.67'
$ws.Range("E49").Value = '  -4.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06008'
$ws.Range("E50").Value = '  -2.06%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '61.60'
$ws.Range("E51").Value = '  -2.99%  '
